$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 16, shifting old row 16 (and its formatting) down to row 17
$ws.Rows.Item(16).Insert()

# New row 16 values (new weekly entry)
$ws.Range("A16").Value = 8
$ws.Range("B16").Value = "Terminal La Palmera de La Serena"
$ws.Range("C16").Value = "Coquimbo"
$ws.Range("D16").Value = 44714
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 100112026
$ws.Range("G16").Value = "Haba"
$ws.Range("H16").Value = "Sin especificar"
$ws.Range("I16").Value = "Primera"
$ws.Range("J16").Value = 400
$ws.Range("K16").Value = 14000
$ws.Range("L16").Value = 15000
$ws.Range("M16").Value = 14500
$ws.Range("N16").Value = '$/saco 25 kilos'
$ws.Range("O16").Value = 'Provincia de Limarí'
$ws.Range("P16").Value = 580
$ws.Range("Q16").Value = 25
$ws.Range("R16").Value = 'Hortaliza'
